# JournalDeBord.xlsx update - TP A1 et A2-4-D du 24/01/2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update date of existing FSIL/CMx2 entry to 23/01/2024 ---
$ws.Range("A2").Value = "01/23/2024"

# --- Row 3: new entry for FSIL/TP on 24/01/2024 (rest of the row content already present) ---
# Copy the date formatting from A2 so A3 reuses the existing date style instead of
# Excel auto-creating a brand new number-format style.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = "01/24/2024"
$ws.Range("B3").Value = "FSIL"

# --- Row 4: new entry for MPAL/TP group A1 on 24/01/2024 ---
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = "01/24/2024"
$ws.Range("B4").Value = "MPAL"
$ws.Range("C4").Value = "TP"
$ws.Range("E4").Value = "x"
$ws.Range("G4").Value = "Simple stack => #2.1 en cours pour la plupart, terminé pour certains. A terminer pour la semaine prochaine."
$ws.Range("H4").Value = "Certains ont du mal à paramétrer le proxy malgré la notice. Peut-être modifier la notice avec Adresse de configuration automatique du proxy : http://cache.iut-rodez.fr/proxy.pac"
$ws.Range("I4").Value = "Inégalité d'autonomie dans l'utilisation de la notice. L'application du override pour le paramétrage Maven a régulièrement été oublié…"

# --- Row 5: new entry for MPAL/TP group A2-4-D on 24/01/2024 ---
$ws.Range("A5").Value = "01/24/2024"
$ws.Range("C5").Value = "TP"
$ws.Range("D5").Value = "x"
$ws.Range("G5").Value = "Simple stack => #2.1 en cours pour la plupart, terminé pour certains. A terminer pour la semaine prochaine."
$ws.Range("I5").Value = "Pour certains, il n'est pas clair que pour la #2.1 il faut implémenter le contenu de la classe SimpleStack avec sa structure de données interne, sans modifier les tests… Je me suis posé la même question lorsque j'ai fait le TP de mon côté, mais j'ai rapidement compris. Peut-être que ce doute nuit à certains. Par ailleurs, savoir qu'il faut utiliser une structure comme ArrayList en interne n'est pas automatique pour tous..."
$ws.Range("H5").Value = 'expliciter "encore plus" ce qui est demandé pour les quelques qui ont du mal à demarrer (coder la structure interne avec un ArrayList) et limiter la charge cognitive ?'

# --- Rows 6, 8, 9, 11, 12, 14, 15, 16: clear the old planned-session placeholders (date + enseignant) ---
$ws.Range("A6").ClearContents()
$ws.Range("B6").ClearContents()

$ws.Range("A8").ClearContents()
$ws.Range("B8").ClearContents()

$ws.Range("A9").ClearContents()
$ws.Range("B9").ClearContents()

$ws.Range("A11").ClearContents()
$ws.Range("B11").ClearContents()

$ws.Range("A12").ClearContents()
$ws.Range("B12").ClearContents()

$ws.Range("A14").ClearContents()
$ws.Range("B14").ClearContents()

$ws.Range("A15").ClearContents()
$ws.Range("B15").ClearContents()

$ws.Range("A16").ClearContents()
$ws.Range("B16").ClearContents()

# --- Row heights for the newly-filled, wrapped-text rows ---
$ws.Rows.Item(4).RowHeight = 63
$ws.Rows.Item(5).RowHeight = 78.75

# --- Selection / view state ---
$ws.Range("H6").Select()
